# Daily attendance processing - 2026-01-18 08:40:54
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Class Statistics block updates (rows 7/8, column L) ---
$ws.Range("L7").Value = 45
$ws.Range("L8").Value = 24

# --- 2) "Recorded By" text order swap: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com" ---
$gRows = @(8,9,10,12,14,15,17,18,23,34,35,36,38,40,41,43,44,49,60,61,62,64,66,67,69,70,75,86,87,88,90,92,93,95,96,101,112,113,114,116,118,119,121,122,127,138,139,140,142,144,145,147,148,153,164,167,170,174,191,194,197,201,218,221,224,228,245,248,251,255,272,275,278,282,299,302,305,309)
foreach ($r in $gRows) {
    $ws.Cells.Item($r, 7).Value = "System, dnasr281@gmail.com"
}

# --- 3) Group Statistics block updates (Missing/Pending counts for B1A1..B1C2) ---
$ws.Range("P15").Value = 4
$ws.Range("Q15").Value = 1
$ws.Range("P16").Value = 3
$ws.Range("Q16").Value = 1
$ws.Range("P17").Value = 3
$ws.Range("Q17").Value = 1
$ws.Range("P18").Value = 3
$ws.Range("Q18").Value = 1
$ws.Range("P19").Value = 3
$ws.Range("Q19").Value = 1
$ws.Range("P20").Value = 4
$ws.Range("Q20").Value = 1

# --- 4) Newly-processed sessions: "Pending" -> "Not Recorded" (rows 26,52,78,104,130,156) ---
# Copy the formatting already used by an existing "Not Recorded" row (row 3) onto
# each processed row (column G, "Recorded By", stays blank on these rows since no
# one has recorded them yet), then set the status text.
$notRecordedRows = @(26,52,78,104,130,156)
foreach ($r in $notRecordedRows) {
    $ws.Range("A3:I3").Copy()
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)
    $ws.Cells.Item($r, 9).Value = "Not Recorded"
}
